# Commit: "Add round icon requirement"
#
# 1) Collapse the spell-check-split "Icons need to be created..." paragraph
#    back into a single run (Find/Replace with identical text forces the
#    runtime to re-flow the run list and drops the now-redundant
#    w:proofErr spell-check markers).
# 2) Insert a new sub-bullet "Both traditional square + new circle" right
#    after the "Application Icon" bullet under the Android heading, to
#    capture the new round-icon requirement.
# 3) Word always keeps the hidden "_GoBack" bookmark at the site of the
#    most recent edit; move it from its old location (inside "Old Record
#    Color") to the end of the text we just typed.

$d = $word.ActiveDocument

# --- 1. Normalize the Iconography intro paragraph -------------------------
$introText = "Icons need to be created for iOS and Android devices, with separate files for each screen density on the respective OSs: normal/@2x/@3x on iOS, mdpi/hdpi/xhdpi/xxhdpi/xxxhdpi on Android. All icons should be .png file types, with transparency where it makes sense"
$d.Content.Find.Execute($introText, $true, $false, $false, $false, $false, $true, 1, $false, $introText, 2) | Out-Null

# --- 2. Insert the new "round icon" bullet ---------------------------------
# Find the "Application Icon" bullet under the Android heading and add a
# new ilvl=1 sub-bullet right after it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Application Icon") {
        $target = $para
        break
    }
}

$target.Range.InsertParagraphAfter()
$newPara = $target.Next()
$newPara.Range.Text = "Both traditional square + new circle"
$newPara.Range.ListFormat.ListLevelNumber = 2

# --- 3. Relocate the "_GoBack" bookmark to the newly typed text -----------
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
# Collapse to just after the last visible character of the new run (i.e.
# before the paragraph mark), matching where Word drops "_GoBack" after
# typing new text.
$endOfNewText = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$endOfNewText.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endOfNewText) | Out-Null
